$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B303").Value = 87504.31089811739
$ws.Range("B304").Value = 88893.57453918592
$ws.Range("B305").Value = 90306.91963708345
$ws.Range("B306").Value = 91744.67394044143
$ws.Range("B307").Value = 93207.1594356229
$ws.Range("B308").Value = 94694.69149521155
$ws.Range("B309").Value = 96207.57799301714
$ws.Range("B310").Value = 97746.11838583663
$ws.Range("B311").Value = 99310.60276239192
$ws.Range("B312").Value = 100901.3108602124
$ws.Range("B313").Value = 102518.5110511117
$ws.Range("B314").Value = 104162.4592961455
$ws.Range("B315").Value = 105833.3980713008
$ws.Range("B316").Value = 107531.55526512
$ws.Range("B317").Value = 109257.1430497207
$ws.Range("B318").Value = 111010.356727047
$ws.Range("B319").Value = 112791.3735520751
$ws.Range("B320").Value = 114600.3515352117
$ws.Range("B321").Value = 116437.4282263228
$ws.Range("B322").Value = 118302.7194827578
$ws.Range("B323").Value = 120196.3182244761
$ws.Range("B324").Value = 122118.2931790851
$ws.Range("B325").Value = 124068.6876204429
$ws.Range("B326").Value = 126047.5181040914
$ws.Range("B327").Value = 128054.7732035716
$ws.Range("B328").Value = 130090.4122516834
$ws.Range("B329").Value = 132154.36409096
$ws.Range("B330").Value = 134246.5258380295
$ws.Range("B331").Value = 136366.7616666898
$ws.Range("B332").Value = 138514.9016147816
$ws.Range("B333").Value = 140690.7404200932
$ws.Range("B334").Value = 142894.0363909452
$ws.Range("B335").Value = 145124.5103170638
$ws.Range("B336").Value = 147381.844426711
$ws.Range("B337").Value = 149665.681396185
$ws.Range("B338").Value = 151975.6234178527
$ws.Range("B339").Value = 154311.2313330686
$ws.Range("B340").Value = 156672.0238366463
$ws.Range("B341").Value = 159057.4767591608
$ws.Range("B342").Value = 161467.0224338428
$ws.Range("B343").Value = 163900.0491547308
$ws.Range("B344").Value = 166355.9007325691
$ws.Range("B345").Value = 168833.8761552153
$ws.Range("B346").Value = 171333.2293587682
$ws.Range("B347").Value = 173853.1691162469
$ws.Range("B348").Value = 176392.8590495294
$ws.Range("B349").Value = 178951.4177711401
$ws.Range("B350").Value = 181527.9191611617
$ws.Range("B351").Value = 184121.3927853182
$ws.Range("B352").Value = 186730.8244590356
$ws.Range("B353").Value = 189355.1569626348
$ws.Range("B354").Value = 191993.2909118975
$ws.Range("B355").Value = 194644.0857882107
$ws.Range("B356").Value = 197306.3611316709
$ws.Range("B357").Value = 199978.8979002218
$ws.Range("B358").Value = 217816.5993111332
$ws.Range("B359").Value = 220706.9739227913
$ws.Range("B360").Value = 223604.2152717551
$ws.Range("B361").Value = 226506.8655326307
$ws.Range("B362").Value = 229413.4354609017
$ws.Range("B363").Value = 232322.4069550165
$ws.Range("B364").Value = 235232.2358005846
$ws.Range("B365").Value = 238141.3545935237
$ws.Range("B366").Value = 241048.1758388956
$ws.Range("B367").Value = 243951.0952212831
$ws.Range("B368").Value = 246848.4950409575
$ws.Range("B369").Value = 249738.7478103951
$ws.Range("B370").Value = 252620.2200035911
$ws.Range("B371").Value = 255491.2759507617
$ws.Range("B372").Value = 258350.2818697478
$ws.Range("B373").Value = 254283.9628051224
$ws.Range("B374").Value = 257039.1086512413
$ws.Range("B375").Value = 259777.8037519292
$ws.Range("B376").Value = 262498.5064773588
$ws.Range("B377").Value = 265199.6983062237
$ws.Range("B378").Value = 267879.8883301202
$ws.Range("B379").Value = 270537.6178119734
$ws.Range("B380").Value = 273171.4647848681
$ws.Range("B381").Value = 275780.048677278
$ws.Range("B382").Value = 278362.0349501816
$ws.Range("B383").Value = 215236.5443044131
$ws.Range("B384").Value = 217171.1826432301
$ws.Range("B385").Value = 219082.6215159229
$ws.Range("B386").Value = 220970.0147998908
$ws.Range("B387").Value = 222832.5691449723
$ws.Range("B388").Value = 224669.5473402231
$ws.Range("B389").Value = 226480.2716090042
$ws.Range("B390").Value = 228264.1268209793
$ws.Range("B391").Value = 230020.5636096019
$ws.Range("B392").Value = 231749.1013841411
$ws.Range("B393").Value = 343367.4969921544
$ws.Range("B394").Value = 345826.1409646957
$ws.Range("B395").Value = 348242.2778140905
$ws.Range("B396").Value = 350615.645944871
$ws.Range("B397").Value = 352946.1041959572
$ws.Range("B398").Value = 355233.6350174868
$ws.Range("B399").Value = 357478.3473894878
$ws.Range("B400").Value = 359680.4794688068
$ws.Range("B401").Value = 361840.4009527031
$ws.Range("B402").Value = 363958.6151459903
$ws.Range("B403").Value = 983043.8436066696
$ws.Range("B404").Value = 988514.1158130482
$ws.Range("B405").Value = 993878.6278342606
$ws.Range("B406").Value = 999140.2052091751
$ws.Range("B407").Value = 1004302.047068411
$ws.Range("B408").Value = 1009367.726283663
$ws.Range("B409").Value = 1014341.188617754
$ws.Range("B410").Value = 1019226.750843105
$ws.Range("B411").Value = 1024029.097804603
$ws.Range("B412").Value = 1028753.278393618
$ws.Range("B413").Value = 1115583.896842943
$ws.Range("B414").Value = 1120532.886744793
$ws.Range("B415").Value = 1125416.141705901
$ws.Range("B416").Value = 1130240.616060636
$ws.Range("B417").Value = 1135013.616674762
$ws.Range("B418").Value = 1139742.790863914
$ws.Range("B419").Value = 1144436.112898908
$ws.Range("B420").Value = 1149101.869059497
$ws.Range("B421").Value = 1153748.641199037
$ws.Range("B422").Value = 1158385.28878004
$ws.Range("B423").Value = 1163020.929344759
$ws.Range("B424").Value = 1167664.917379638
$ws.Range("B425").Value = 1172326.821540681
$ws.Range("B426").Value = 1177016.400200688
$ws.Range("B427").Value = 1181743.575288674
$ws.Range("B428").Value = 1186518.404391232
$ws.Range("B429").Value = 1191351.05108934
$ws.Range("B430").Value = 1196251.753512539
$ws.Range("B431").Value = 1201230.791095426
$ws.Range("B432").Value = 1206298.449529908
$ws.Range("B433").Value = 1211464.983917283
$ws.Range("B434").Value = 1216740.580130534
$ws.Range("B435").Value = 1222135.314413717
$ws.Range("B436").Value = 1227659.111252251
$ws.Range("B437").Value = 1233321.699569666
$ws.Range("B438").Value = 1239132.567314676
$ws.Range("B439").Value = 1245100.914527719
$ws.Range("B440").Value = 1251235.60498753
$ws.Range("B441").Value = 1257545.116565706
$ws.Range("B442").Value = 1264037.490434306
$ws.Range("B443").Value = 1270720.279295618
$ws.Range("B444").Value = 1277600.494829536
$ws.Range("B445").Value = 1284684.554576099
$ws.Range("B446").Value = 1291978.228498197
$ws.Range("B447").Value = 1299486.585494871
$ws.Range("B448").Value = 1307213.94016231
$ws.Range("B449").Value = 1315163.800127917
$ws.Range("B450").Value = 1323338.814301805
$ws.Range("B451").Value = 1331740.722425331
$ws.Range("B452").Value = 1340370.306310536
